$d = $word.ActiveDocument
$rng = $d.Content

$xml = @'
<pkg:package xmlns:pkg="http://schemas.microsoft.com/office/2006/xmlPackage"><pkg:part pkg:name="/word/document.xml" pkg:contentType="application/vnd.openxmlformats-officedocument.wordprocessingml.document.main+xml"><pkg:xmlData><w:document xmlns:w="http://schemas.openxmlformats.org/wordprocessingml/2006/main"><w:body>
    <w:p>
      <w:pPr>
        <w:pStyle w:val="NormalWeb"/>
        <w:spacing w:after="0" w:afterAutospacing="0"/>
        <w:contextualSpacing/>
      </w:pPr>
      <w:r>
        <w:t>\section{Geographic Visualization}</w:t>
      </w:r>
    </w:p>
    <w:p>
      <w:pPr>
        <w:pStyle w:val="NormalWeb"/>
        <w:spacing w:after="0" w:afterAutospacing="0"/>
        <w:contextualSpacing/>
      </w:pPr>
      <w:r>
        <w:t xml:space="preserve">Geovisualization is deeply rooted in traditional cartography which has been around for thousands of years in human history. Geovisualization is the visualization of geospatial information to create human understanding that leads to data exploration and decision making [1]. </w:t>
      </w:r>
    </w:p>
    <w:p>
      <w:pPr>
        <w:pStyle w:val="NormalWeb"/>
        <w:spacing w:after="0" w:afterAutospacing="0"/>
        <w:contextualSpacing/>
      </w:pPr>
      <w:r>
        <w:t xml:space="preserve">Traditional maps known as static maps are simply fixed images. These maps can be produced on traditional mediums such as hardcopy, like books, atlases and magazines. They can also appear online as images on websites in file formats such as Portable Network Graphic (PNG), Joint Photographic Experts Group (JPEG), and Portable Document Format (PDF). </w:t>
      </w:r>
    </w:p>
    <w:p>
      <w:pPr>
        <w:pStyle w:val="NormalWeb"/>
        <w:spacing w:after="0" w:afterAutospacing="0"/>
        <w:contextualSpacing/>
      </w:pPr>
      <w:r>
        <w:t>\begin{figure}</w:t>
      </w:r>
    </w:p>
    <w:p>
      <w:pPr>
        <w:pStyle w:val="NormalWeb"/>
        <w:spacing w:after="0" w:afterAutospacing="0"/>
        <w:contextualSpacing/>
      </w:pPr>
      <w:r>
        <w:t>\includegraphics[width=0.85\textwidth]{geoVizExample.png}</w:t>
      </w:r>
    </w:p>
    <w:p>
      <w:pPr>
        <w:pStyle w:val="NormalWeb"/>
        <w:spacing w:after="0" w:afterAutospacing="0"/>
        <w:contextualSpacing/>
      </w:pPr>
      <w:r>
        <w:t>\caption{This map from The U.S. Census Bureau is an example of a geovisualization on a static map. It shows the distribution of urban areas in the U.S. and Puerto Rico.}</w:t>
      </w:r>
    </w:p>
    <w:p>
      <w:pPr>
        <w:pStyle w:val="NormalWeb"/>
        <w:spacing w:after="0" w:afterAutospacing="0"/>
        <w:contextualSpacing/>
      </w:pPr>
      <w:r>
        <w:t>\label{fig.1}</w:t>
      </w:r>
    </w:p>
    <w:p>
      <w:pPr>
        <w:pStyle w:val="NormalWeb"/>
        <w:spacing w:after="0" w:afterAutospacing="0"/>
        <w:contextualSpacing/>
      </w:pPr>
      <w:r>
        <w:t>\end{figure}</w:t>
      </w:r>
    </w:p>
    <w:p>
      <w:pPr>
        <w:pStyle w:val="NormalWeb"/>
        <w:spacing w:after="0" w:afterAutospacing="0"/>
        <w:contextualSpacing/>
      </w:pPr>
      <w:r>
        <w:t xml:space="preserve">\par Interactive maps allow the user to zoom in and out, hover-over popups, and more to engage data and find underlying patterns in greater depth. Interactive maps are produced and viewed on computers. </w:t>
      </w:r>
    </w:p>
    <w:p>
      <w:pPr>
        <w:pStyle w:val="NormalWeb"/>
        <w:spacing w:after="0" w:afterAutospacing="0"/>
        <w:contextualSpacing/>
      </w:pPr>
      <w:r>
        <w:t>\begin{figure}</w:t>
      </w:r>
    </w:p>
    <w:p>
      <w:pPr>
        <w:pStyle w:val="NormalWeb"/>
        <w:spacing w:after="0" w:afterAutospacing="0"/>
        <w:contextualSpacing/>
      </w:pPr>
      <w:r>
        <w:t>\includegraphics[width=0.85\textwidth]{geoVizExample2.png}</w:t>
      </w:r>
    </w:p>
    <w:p>
      <w:pPr>
        <w:pStyle w:val="NormalWeb"/>
        <w:spacing w:after="0" w:afterAutospacing="0"/>
        <w:contextualSpacing/>
      </w:pPr>
      <w:r>
        <w:t>\caption{This map from The U.S. Geoglogical Survey is an example of an interactive geovisualization on the Internet. The interactive map allows the user to look at the Active Groundwater Levels at wells across the U.S..}</w:t>
      </w:r>
    </w:p>
    <w:p>
      <w:pPr>
        <w:pStyle w:val="NormalWeb"/>
        <w:spacing w:after="0" w:afterAutospacing="0"/>
        <w:contextualSpacing/>
      </w:pPr>
      <w:r>
        <w:t>\label{fig.2}</w:t>
      </w:r>
    </w:p>
    <w:p>
      <w:pPr>
        <w:pStyle w:val="NormalWeb"/>
        <w:spacing w:after="0" w:afterAutospacing="0"/>
        <w:contextualSpacing/>
      </w:pPr>
      <w:r>
        <w:t>\end{figure}</w:t>
      </w:r>
    </w:p>
    <w:p>
      <w:pPr>
        <w:pStyle w:val="NormalWeb"/>
        <w:spacing w:after="0" w:afterAutospacing="0"/>
        <w:contextualSpacing/>
      </w:pPr>
      <w:r>
        <w:t xml:space="preserve">\par Some early work in geovisualization can be traced back to the term geographic visualization by the National Science Foundation in 1987. If we go back a decade earlier, we find that Jacque Bertin presented design principals for presenting cartographic and information design to explore data [2]. The International Cartographic Association (ICA) created a Commission on Visualization in 1995 to stimulate geovisualization research and encourage interdisciplinary research to create highly interactive, exploratory methods to initiate knowledge construction. In 2015, the ICA created the Commission on Visual Analytics to support geovisualization. This commission focuses on interactive visualizations that can support knowledge construction and insights from spatial data in forms that are both big and small\footnote{Robinson, Anthony, "New Directions in Geovisual Analytics: Visualization, Computation, and Evaluation", International Cartographic Association Commission on Visual Analytics, 2/2/2018, https://viz.icaci.org/}.</w:t>
      </w:r>
    </w:p>
    <w:p>
      <w:pPr>
        <w:pStyle w:val="NormalWeb"/>
        <w:spacing w:before="0" w:beforeAutospacing="0" w:after="0" w:afterAutospacing="0"/>
        <w:contextualSpacing/>
      </w:pPr>
      <w:r>
        <w:t>\par Modern information availability has helped lead to an explosion of geovisualization tools. There are numerous private companies using geovisualization as revenue drivers. Perhaps the best known is Google's Maps Platform. Google's Map Platform has over one billion monthly active users and gets 25 million updates a day\footnote{Google Maps Platform, https://cloud.google.com/maps-platform/maps/}. Google sells their platform to companies to display information spatially. Government agencies use geovisualization for military, forestry, fishery, demographic and economic data display for both internal and external communication.</w:t>
      </w:r>
      <w:bookmarkStart w:id="0" w:name="_GoBack"/>
      <w:bookmarkEnd w:id="0"/>
    </w:p>
</w:body></w:document></pkg:xmlData></pkg:part></pkg:package>
'@

$rng.InsertXML($xml)
